$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14,1).Value = "Hydraulic and electric properties of tissues"
$ws.Cells.Item(14,2).Value = "Contributed talk"
$ws.Cells.Item(14,3).Value = "talk-11"
$ws.Cells.Item(14,4).Value = "Physics meets Biology"
$ws.Cells.Item(14,5).Value = 45223
$ws.Cells.Item(14,6).Value = "Rice Global Paris Center, Paris, France"

[void]$ws.Range("D17").Select()

$r14e = $ws.Cells.Item(14,5)
"E14 value2: $($r14e.Value2)"
$r14a = $ws.Cells.Item(14,1)
"A14 value2: $($r14a.Value2)"
